# "Updated figs and tabs"
# - Duplicate the existing "by_prov" sheet, keep the old data on a sheet
#   renamed "by_prov_old", and turn the (now first) duplicate into the new
#   "by_prov" sheet carrying refreshed (weighted) proportions + CIs.
# - Fix two province labels ("NFL & NS" -> "NL & NS", "MT & AB" -> "MB & AB")
#   on the new sheet.
# - Apply a 2-decimal number format to the new numeric data.

$wb = $excel.ActiveWorkbook

$wsOld = $wb.Worksheets.Item("by_prov")
$wsOld.Copy($wsOld)

# After Copy, the freshly inserted duplicate sits at position 1 and the
# original (still holding the legacy data) has been pushed to position 2.
$wsNew = $wb.Worksheets.Item(1)
$wsOrig = $wb.Worksheets.Item(2)

$wsOrig.Name = "by_prov_old"
$wsNew.Name = "by_prov"

# Refresh the data on the new "by_prov" sheet with weighted proportions and
# confidence intervals (replacing the old OR/lower/higher numbers).
# Cells are written in the same order Excel's shared-string table ends up
# using: "MB & AB" (row 8/9) registers before "NL & NS" (row 2/3).

$wsNew.Range("E8").Value = "MB & AB"
$wsNew.Range("E9").Value = "MB & AB"
$wsNew.Range("E2").Value = "NL & NS"
$wsNew.Range("E3").Value = "NL & NS"

$wsNew.Range("B2:D11").NumberFormat = "0.00_ "

$wsNew.Range("B2").Value = 0.31948881789137379
$wsNew.Range("C2").Value = 0.1855287569573284
$wsNew.Range("D2").Value = 0.52631578947368418

$wsNew.Range("B3").Value = 0.7142857142857143
$wsNew.Range("C3").Value = 0.60606060606060608
$wsNew.Range("D3").Value = 0.84745762711864414

$wsNew.Range("B4").Value = 0.30674846625766872
$wsNew.Range("C4").Value = 0.19960079840319361
$wsNew.Range("D4").Value = 0.46296296296296291

$wsNew.Range("B5").Value = 0.69444444444444442
$wsNew.Range("C5").Value = 0.58823529411764708
$wsNew.Range("D5").Value = 0.81967213114754101

$wsNew.Range("B6").Value = 0.2386634844868735
$wsNew.Range("C6").Value = 0.16260162601626016
$wsNew.Range("D6").Value = 0.34129692832764502

$wsNew.Range("B7").Value = 0.69444444444444442
$wsNew.Range("C7").Value = 0.5988023952095809
$wsNew.Range("D7").Value = 0.80645161290322587

$wsNew.Range("B8").Value = 0.52910052910052918
$wsNew.Range("C8").Value = 0.37037037037037035
$wsNew.Range("D8").Value = 0.74626865671641784

$wsNew.Range("B9").Value = 0.67114093959731547
$wsNew.Range("C9").Value = 0.5714285714285714
$wsNew.Range("D9").Value = 0.78740157480314954

$wsNew.Range("B10").Value = 0.13774104683195593
$wsNew.Range("C10").Value = 0.080645161290322578
$wsNew.Range("D10").Value = 0.22123893805309736

$wsNew.Range("B11").Value = 0.68493150684931503
$wsNew.Range("C11").Value = 0.58479532163742687
$wsNew.Range("D11").Value = 0.8

$wsNew.Range("E10").Select()
